$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Petal coordinate refinement (small numeric re-measurement) ---
$ws.Range("A2").Value = 135.01820009346088
$ws.Range("B2").Value = -472.47263505383114
$ws.Range("C2").Value = -117.99902208507902
$ws.Range("D2").Value = 697.58478850862195
$ws.Range("E2").Value = -684.22727617799887
$ws.Range("F2").Value = -118.11576838924343

# --- Row 22/23: R3 module coordinate refinement ---
$ws.Range("A22").Value = 316.67161377698926
$ws.Range("B22").Value = -163.82848220086203
$ws.Range("C22").Value = -125.20785120469837

$ws.Range("A23").Value = 429.81037436784902
$ws.Range("B23").Value = -237.84213693251039
$ws.Range("C23").Value = -125.7655754805098

# --- Row 40: R2 pick-up tool holder angle ---
$ws.Range("A40").Value = -90

# --- Row 56: R5 pick-up tool holder angle ---
$ws.Range("A56").Value = -90

# --- Row 59: camera offset relative to gantry axis (now formula-driven) ---
$ws.Range("A59").Formula = "=-0.549865"
$ws.Range("B59").Formula = "=-101.798218"
$ws.Range("C59").Value = -23.54
$ws.Range("D59").Value = "<-- x, y, z offset of camera relative to gantry axis"

# --- Row 64-66: Cognex camera network / calibration settings ---
$ws.Range("A64").Value = "169.254.41.3"
$ws.Range("A65").Value = 50290
$ws.Range("A66").Value = 0.365022

# --- Selection moved to A23:C23 (as last user action before save) ---
$ws.Range("A23:C23").Select()
